$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for rows 2-15 from serial date 45175 (2023-09-06)
# to serial date 45183 (2023-09-14).
foreach ($row in 2..15) {
    $ws.Cells.Item($row, 3).Value = 45183
}
